# Inclusão de lógica para inserção do cabeçalho das planilhas
#
# - Remove the "idade" worksheet (its data is folded into the "nome" sheet
#   as a header row instead).
# - Add the header row ("Nomes" / "Idade") to the "nome" worksheet.

$wb = $excel.ActiveWorkbook

# Avoid the "are you sure you want to delete this sheet" confirmation dialog.
$excel.DisplayAlerts = $false

$idadeSheet = $wb.Worksheets.Item("idade")
if ($idadeSheet -ne $null) {
    $idadeSheet.Delete()
}

$nomeSheet = $wb.Worksheets.Item("nome")
$nomeSheet.Range("A1").Value = "Nomes"
$nomeSheet.Range("B1").Value = "Idade"

$excel.DisplayAlerts = $true
